$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values could be misinterpreted as numbers by Excel;
# force them to remain plain text and keep the default (unstyled) cell style.
$priceUpdates = @{
    "D2" = "29.934.39"
    "D3" = "1.892.71"
    "D4" = "1.001"
    "D5" = "0.7734"
    "D6" = "243.82"
    "D8" = "0.3134"
    "D9" = "25.64"
    "D10" = "0.07329"
    "D11" = "0.08057"
    "D12" = "0.7721"
    "D13" = "5.496"
    "D14" = "94.01"
    "D15" = "1.796.43"
    "D16" = "6.221"
    "D17" = "29.834.14"
    "D18" = "14.00"
    "D19" = "246.40"
    "D20" = "0.000007857"
    "D21" = "8.163"
    "D22" = "0.9997"
    "D23" = "2.085.97"
    "D24" = "1.001"
    "D25" = "0.1570"
    "D26" = "9.443"
    "D27" = "162.34"
    "D28" = "18.76"
    "D29" = "2.024"
    "D31" = "1.543"
    "D32" = "4.475"
    "D33" = "0.05554"
    "D34" = "4.066"
    "D35" = "1.235"
    "D36" = "0.7495"
    "D37" = "0.9994"
    "D38" = "2.681"
    "D39" = "0.01928"
    "D41" = "0.4473"
    "D42" = "74.16"
    "D43" = "1.099.26"
    "D44" = "6.002"
    "D45" = "0.8500"
    "D47" = "1.886"
    "D49" = "7.540"
    "D50" = "9.790"
    "D51" = "2.990"
}

foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.Style = "Normal"
}

# Remaining textual updates (volume %, coin names/links swapped between rows 14 and 15).
$textUpdates = @{
    "E2" = "  +0.39%  "
    "E3" = "  -0.02%  "
    "E4" = "  +0.01%  "
    "E5" = "  -2.75%  "
    "E6" = "  +0.38%  "
    "E7" = "  -0.02%  "
    "E8" = "  -0.99%  "
    "E9" = "  +0.82%  "
    "E10" = "  +3.81%  "
    "E11" = "  -0.29%  "
    "E12" = "  +0.56%  "
    "E13" = "  +2.67%  "
    "B14" = "Litecoin"
    "C14" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "E14" = "  +1.61%  "
    "B15" = "WrappedEther"
    "C15" = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
    "E15" = "  -5.46%  "
    "E16" = "  +3.56%  "
    "E17" = "  +0.01%  "
    "E18" = "  +0.93%  "
    "E19" = "  +0.78%  "
    "E20" = "  +1.90%  "
    "E21" = "  -2.30%  "
    "E22" = "  -0.05%  "
    "E23" = "  -2.83%  "
    "E24" = "  -0.01%  "
    "E25" = "  -4.22%  "
    "E26" = "  +0.90%  "
    "E27" = "  -2.32%  "
    "E28" = "  +0.31%  "
    "E29" = "  -1.57%  "
    "E30" = "  +1.79%  "
    "E31" = "  +0.20%  "
    "E32" = "  +0.83%  "
    "E33" = "  -2.72%  "
    "E34" = "  +0.50%  "
    "E35" = "  -2.06%  "
    "E36" = "  +1.42%  "
    "E37" = "  +0.09%  "
    "E38" = "  +1.96%  "
    "E39" = "  +0.97%  "
    "E40" = "  +0.17%  "
    "E41" = "  +1.50%  "
    "E42" = "  +2.29%  "
    "E43" = "  +6.25%  "
    "E44" = "  +3.24%  "
    "E45" = "  +1.03%  "
    "E46" = "  -0.03%  "
    "E47" = "  +0.64%  "
    "E48" = "  -0.70%  "
    "E49" = "  +1.50%  "
    "E50" = "  -1.94%  "
    "E51" = "  +3.07%  "
}

foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

